$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet's data rows (2..12) represent weekly price records, newest on
# top. A new week of data is being added, so the existing rows need to
# shift down by two positions (2..12 -> 4..14) before writing the two new
# rows (2 and 3). Working from the bottom up avoids overwriting data
# before it has been copied, and plain cell-value writes (rather than
# Range.Insert) keep the existing per-cell formatting (e.g. the date
# style already present in column D of every destination row) untouched.
$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
for ($srcRow = 12; $srcRow -ge 2; $srcRow--) {
    $dstRow = $srcRow + 2
    foreach ($col in $columns) {
        $ws.Range($col + $dstRow).Value = $ws.Range($col + $srcRow).Formula
    }
}

# Row 2: new weekly record (Primera)
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "Femacal de La Calera"
$ws.Range("C2").Value = "Coquimbo"
$ws.Range("D2").Value = 44756
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 100112043
$ws.Range("G2").Value = "Pepino dulce"
$ws.Range("H2").Value = "Cultivar IV Región"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 65
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 14000
$ws.Range("N2").Value = "`$/caja 15 kilos"
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 933
$ws.Range("Q2").Value = 15
$ws.Range("R2").Value = "Hortaliza"

# Row 3: new weekly record (Segunda)
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Femacal de La Calera"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44756
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 100112043
$ws.Range("G3").Value = "Pepino dulce"
$ws.Range("H3").Value = "Cultivar IV Región"
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 68
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 12000
$ws.Range("N3").Value = "`$/caja 15 kilos"
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 800
$ws.Range("Q3").Value = 15
$ws.Range("R3").Value = "Hortaliza"

$ws.Range("A1").Select()
